# Update "想去人数" (want-to-go count) values in column F on the
# "展览" and "全部类型" worksheets, as produced by the latest site scrape.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 167
$ws1.Range("F5").Value  = 1821
$ws1.Range("F9").Value  = 2366
$ws1.Range("F13").Value = 1434
$ws1.Range("F15").Value = 35
$ws1.Range("F22").Value = 211
$ws1.Range("F24").Value = 93
$ws1.Range("F26").Value = 1479
$ws1.Range("F29").Value = 233
$ws1.Range("F32").Value = 371

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 167
$ws4.Range("F5").Value  = 1821
$ws4.Range("F10").Value = 2366
$ws4.Range("F14").Value = 1434
$ws4.Range("F16").Value = 35
$ws4.Range("F23").Value = 211
$ws4.Range("F25").Value = 93
$ws4.Range("F27").Value = 1479
$ws4.Range("F30").Value = 233
$ws4.Range("F33").Value = 371
